# Add a new column E ("exception") to the "Formula Tests" sheet test table,
# with a header in E5 and a "true" marker value in E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Formula Tests")

# E6 = "true" -- must land as a plain shared string (not an Excel Boolean),
# so build it via a text formula and paste back only the resulting value.
$cell6 = $ws.Cells.Item(6, 5)
$cell6.Formula = '="true"'
$cell6.Copy()
$cell6.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# E5 = "exception" header, styled like the rest of the header row (D5).
$cell5 = $ws.Cells.Item(5, 5)
$cell5.Formula = '="exception"'
$cell5.Copy()
$cell5.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Cells.Item(5, 4).Copy()
$cell5.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Move/keep the active selection on E8, matching the saved view state.
$ws.Range("E8").Select() | Out-Null
